$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "50.938.42"
$ws.Range("E2").Value = "  -1.57%  "

# Row 3
$ws.Range("D3").Value = "2.923.69"
$ws.Range("E3").Value = "  -2.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "378.85"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.87"
$ws.Range("E6").Value = "  -3.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  -2.00%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -3.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.38"
$ws.Range("E10").Value = "  -3.75%  "

# Row 11
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0834"
$ws.Range("E12").Value = "  -1.48%  "

# Row 13
$ws.Range("D13").Value = "3.389.52"
$ws.Range("E13").Value = "  -2.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.93"
$ws.Range("E14").Value = "  -5.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.32"
$ws.Range("E15").Value = "  -3.30%  "

# Row 16
$ws.Range("D16").Value = "2.897.47"
$ws.Range("E16").Value = "  -3.88%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.970"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").Value = "50.895.10"
$ws.Range("E18").Value = "  -1.67%  "

# Row 19
$ws.Range("E19").Value = "  -9.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  -5.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  -5.52%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("E22").Value = "  -1.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.16"
$ws.Range("E23").Value = "  -1.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.30"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
$ws.Range("E25").Value = "  +2.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.15"
$ws.Range("E26").Value = "  +9.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.58"
$ws.Range("E27").Value = "  +3.01%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +8.26%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.166"
$ws.Range("E30").Value = "  -4.51%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.50"
$ws.Range("E31").Value = "  -2.71%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.75"
$ws.Range("E32").Value = "  -2.47%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.53"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.86"
$ws.Range("E34").Value = "  -2.30%  "

# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.04"
$ws.Range("E35").Value = "  -2.02%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0448"
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.95"
$ws.Range("E38").Value = "  -5.62%  "

# Row 39
$ws.Range("E39").Value = "  -3.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.47"
$ws.Range("E40").Value = "  -6.04%  "

# Row 41
$ws.Range("E41").Value = "  -1.93%  "

# Row 42
$ws.Range("E42").Value = "  -5.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.33"
$ws.Range("E43").Value = "  -3.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.18"
$ws.Range("E44").Value = "  -5.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.04"
$ws.Range("E45").Value = "  -2.44%  "

# Row 46
$ws.Range("E46").Value = "  -3.34%  "

# Row 47
$ws.Range("E47").Value = "  -3.78%  "

# Row 48
$ws.Range("D48").Value = "1.998.00"
$ws.Range("E48").Value = "  -2.63%  "

# Row 49
$ws.Range("E49").Value = "  -3.46%  "

# Row 50
$ws.Range("E50").Value = "  -2.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.477"
$ws.Range("E51").Value = "  +9.94%  "
